$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the confusion-matrix label: cell B3 was mislabeled as a
# "FALSE positives" entry when it is actually the FALSE negatives cell.
$ws.Range("B3").Value = "4 FALSE negatives (Type II error)"

# Leave the active selection on the corrected cell.
$ws.Range("B3").Select()
